# Adds two extra biomarkers (CD8 and Napsin A) to biomarker_rules_enoc:
#   - "translation" sheet: new original_score -> translated_score rows
#   - "consolidation" sheet: new rule_type/rule_value -> consolidated_value rows
# Also tweaks a couple of pre-existing "consolidation" rows (CTNNB1 "else" rule
# and the p16 block) to make room / fix them up as part of the same commit.

$wb = $excel.ActiveWorkbook
$wsTranslation   = $wb.Worksheets.Item("translation")
$wsConsolidation = $wb.Worksheets.Item("consolidation")

# ---- helpers -------------------------------------------------------------

# Force a cell to be stored as TEXT (shared string) even when the text looks
# like a number (e.g. "2"), without leaving the column's normal number
# format behind.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = "#,##0"
}

# Set a "translation" sheet row: biomarker | original_score | translated_score
function Set-TranslationRow($row, $biomarker, $origScore, $translatedScore) {
    $ws = $wsTranslation
    $ws.Cells.Item($row, 1).Value = $biomarker
    if ($origScore -is [string]) {
        Set-TextValue $ws.Cells.Item($row, 2) $origScore
    } else {
        $ws.Cells.Item($row, 2).Value = $origScore
    }
    $ws.Cells.Item($row, 3).Value = $translatedScore
}

# Set a "consolidation" sheet row: biomarker | rule_type | rule_value | consolidated_value
function Set-ConsolidationRow($row, $biomarker, $ruleType, $ruleValue, $consolidatedValue) {
    $ws = $wsConsolidation
    $ws.Cells.Item($row, 1).Value = $biomarker
    if ($ruleType -is [string]) {
        Set-TextValue $ws.Cells.Item($row, 2) $ruleType
    } else {
        $ws.Cells.Item($row, 2).Value = $ruleType
    }
    if ($null -eq $ruleValue) {
        $ws.Cells.Item($row, 3).Value = ""
    } else {
        $ws.Cells.Item($row, 3).Value = $ruleValue
    }
    $ws.Cells.Item($row, 4).Value = $consolidatedValue
}

# ---- "translation" sheet: append CD8 and Napsin A rules (rows 125-135) ---

Set-TranslationRow 125 "CD8 "     0   "negative"
Set-TranslationRow 126 "CD8 "     1   "1 or 2 IEL"
Set-TranslationRow 127 "CD8"      2   "3 - 19 IEL "
Set-TranslationRow 128 "CD8"      3   "20 or more IEL "
Set-TranslationRow 129 "CD8"      9   "Unk "
Set-TranslationRow 130 "CD8"      "x" "Unk"
Set-TranslationRow 131 "Napsin A" 0   "negative "
Set-TranslationRow 132 "Napsin A" 1   "focal (1-50%)"
Set-TranslationRow 133 "Napsin A" 2   "diffuse (>50%)"
Set-TranslationRow 134 "Napsin A" 9   "Unk "
Set-TranslationRow 135 "Napsin A" "x" "Unk "

# ---- "consolidation" sheet -------------------------------------------------

# Row 53 (CTNNB1) changes from "else" -> "negative " to an explicit "any" rule
# matching literal "negative ", and a brand-new CTNNB1 "else" rule is added
# right after it as row 54 (pushing the old p16 block down by one row).
Set-ConsolidationRow 53 "CTNNB1" "any"   "negative " "negative "
Set-ConsolidationRow 54 "CTNNB1" "else " $null       "__check__"

# p16 block, shifted down by one row from its original 54-57 position, plus
# one brand new "else" rule appended at the end (row 59).
Set-ConsolidationRow 55 "p16 " "any "  "normal "                                      "normal "
Set-ConsolidationRow 56 "p16 " 2       "Abnormal complete absence & abnormal block "  "duo abnormal "
Set-ConsolidationRow 57 "p16 " "any "  "abnormal block "                              "abnormal block "
Set-ConsolidationRow 58 "p16 " "any"   "abnormal complete absence "                   "abnormal complete absence "
Set-ConsolidationRow 59 "p16 " "else"  $null                                          "__check__"

# New CD8 rules (rows 60-64)
Set-ConsolidationRow 60 "CD8 " "any"   "20 or more IEL " "high "
Set-ConsolidationRow 61 "CD8 " "any"   "3 - 19 IEL "     "moderate "
Set-ConsolidationRow 62 "CD8"  "any "  "1 or 2 IEL"      "none/low"
Set-ConsolidationRow 63 "CD8"  "any "  "negative "       "none/low"
Set-ConsolidationRow 64 "CD8"  "else " $null             "__check__"

# New Napsin A rules (rows 65-68)
Set-ConsolidationRow 65 "Napsin A" "any "  "diffuse (>50%)" "positive"
Set-ConsolidationRow 66 "Napsin A" "any "  "focal (1-50%)"  "positive"
Set-ConsolidationRow 67 "Napsin A" "any"   "negative "      "negative "
Set-ConsolidationRow 68 "Napsin A" "else " $null            "__check__"

# ---- restore view/selection state -----------------------------------------

$wsTranslation.Range("A77").Select() | Out-Null
$wsConsolidation.Activate() | Out-Null
$wsConsolidation.Range("G54").Select() | Out-Null

Write-Output "biomarker rules updated: CD8 and Napsin A added"
